$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.067.62"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "1.651.99"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.40"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5281"
$ws.Range("E6").Value = "  +1.25%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2600"
$ws.Range("E8").Value = "  -1.39%  "

$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.517"
$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("D13").Value = "1.649.98"
$ws.Range("E13").Value = "  -2.01%  "

$ws.Range("D14").Value = "1.879.05"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5478"
$ws.Range("E15").Value = "  +0.60%  "

$ws.Range("D16").Value = "0.0₅8193"
$ws.Range("E16").Value = "  +0.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.31"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").Value = "26.077.08"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.580"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.64"
$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("E22").Value = "  +0.59%  "

$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.87"
$ws.Range("E25").Value = "  +3.72%  "

$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.209"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.99"
$ws.Range("E28").Value = "  -1.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.454"
$ws.Range("E29").Value = "  +3.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05794"
$ws.Range("E30").Value = "  -2.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.270"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.543"
$ws.Range("E32").Value = "  +0.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.262"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.595"
$ws.Range("E34").Value = "  +1.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.794"
$ws.Range("E35").Value = "  +0.85%  "

$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9427"
$ws.Range("E37").Value = "  -1.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5741"
$ws.Range("E38").Value = "  +1.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01610"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8483"
$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.27"
$ws.Range("E41").Value = "  +3.73%  "

$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.709"
$ws.Range("E43").Value = "  -4.15%  "

$ws.Range("D44").Value = "1.029.81"
$ws.Range("E44").Value = "  +2.61%  "

$ws.Range("D45").Value = "1.793.58"
$ws.Range("E45").Value = "  -0.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.91"
$ws.Range("E46").Value = "  +0.59%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4327"
$ws.Range("E48").Value = "  -0.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.828"
$ws.Range("E49").Value = "  -1.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05141"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.445"
$ws.Range("E51").Value = "  -1.24%  "
